$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

$ws.Range("A7").Value = 10122014
$ws.Range("B7").Value = "Items/Prefabs/Weapons/Skill/Wand_1_Epic_Skill"

$ws.Range("B8").Select()
